$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting existing rows 58-62 down to 59-63.
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with the new weekly record.
$ws.Range("A58").Value = 5
$ws.Range("B58").Value = "Macroferia Regional de Talca"
$ws.Range("C58").Value = "Maule"
$ws.Range("D58").Value = 44516
$ws.Range("E58").Value = 7
$ws.Range("F58").Value = 100112022
$ws.Range("G58").Value = "Arveja Verde"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 600
$ws.Range("K58").Value = 15000
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = 15000
$ws.Range("N58").Value = "$/saco 25 kilos"
$ws.Range("O58").Value = "Región del Maule"
$ws.Range("P58").Value = 600
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"
